$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update "Status as of ..." header text ---
$ws.Range("X1").Value = "Status as of July 11, 2025"

# --- Column widths (chars) ---
$ws.Columns.Item(1).ColumnWidth = 12.166666666666666
$ws.Columns.Item(2).ColumnWidth = 19.166666666666668
$ws.Columns.Item(3).ColumnWidth = 10.166666666666666
$ws.Columns.Item(4).ColumnWidth = 44.166666666666664
$ws.Columns.Item(5).ColumnWidth = 13.166666666666666
$ws.Columns.Item(6).ColumnWidth = 13.166666666666666
$ws.Columns.Item(7).ColumnWidth = 13.166666666666666
$ws.Columns.Item(8).ColumnWidth = 103.16666666666667
$ws.Columns.Item(9).ColumnWidth = 17.166666666666668
$ws.Columns.Item(10).ColumnWidth = 16.166666666666668
$ws.Columns.Item(11).ColumnWidth = 10.166666666666666
$ws.Columns.Item(12).ColumnWidth = 25.166666666666668
$ws.Columns.Item(13).ColumnWidth = 25.166666666666668
$ws.Columns.Item(14).ColumnWidth = 26.166666666666668
$ws.Columns.Item(15).ColumnWidth = 21.166666666666668
$ws.Columns.Item(16).ColumnWidth = 12.166666666666666
$ws.Columns.Item(17).ColumnWidth = 30.166666666666668
$ws.Columns.Item(18).ColumnWidth = 26.166666666666668
$ws.Columns.Item(19).ColumnWidth = 14.166666666666666
$ws.Columns.Item(20).ColumnWidth = 32.166666666666664
$ws.Columns.Item(21).ColumnWidth = 30.166666666666668
$ws.Columns.Item(22).ColumnWidth = 38.166666666666664
$ws.Columns.Item(23).ColumnWidth = 106.16666666666667
$ws.Columns.Item(24).ColumnWidth = 27.166666666666668

# --- Freeze the header row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# --- Header row formatting: center/middle + wrap text ---
$x1 = $ws.Range("X1")
$x1.HorizontalAlignment = -4108
$x1.VerticalAlignment = -4108

$headerRest = $ws.Range("A1:W1")
$headerRest.VerticalAlignment = -4108

$headerAll = $ws.Range("A1:X1")
$headerAll.WrapText = $true

# Highlight the "Status as of ..." column (yellow fill)
$ws.Range("X1:X2").Interior.ColorIndex = 6
